$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 42
$ws1.Range("F7").Value = 43
$ws1.Range("F18").Value = 5191
$ws1.Range("F22").Value = 2304
$ws1.Range("F25").Value = 2147

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 42
$ws4.Range("F7").Value = 43
$ws4.Range("F18").Value = 5191
$ws4.Range("F24").Value = 2304
$ws4.Range("F28").Value = 2147
